# Edit script for LOUISIANA_2018.xlsx
# 1) Rename header row columns to short machine-friendly names
# 2) Title-case Spanish stopwords (de/del/las/el/los/la/y) in municipality
#    (and a couple of state) names throughout the data rows
# 3) Remove the trailing metadata/footer rows (776-780), shrinking the
#    sheet's used range down to A1:D774

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames -----------------------------------------------
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- 2) Title-case the "de/del/las/el/los/la/y" connector words ------
$updates = @(
    @(5, 2, 'Pabellón De Arteaga'),
    @(6, 2, 'Rincón De Romos'),
    @(22, 2, 'Bejucal De Ocampo'),
    @(23, 2, 'Benemérito De Las Américas'),
    @(28, 2, 'Chiapa De Corzo'),
    @(31, 2, 'Comitán De Domínguez'),
    @(41, 2, 'Marqués De Comillas'),
    @(45, 2, 'Ocozocoautla De Espinosa'),
    @(52, 2, 'Salto De Agua'),
    @(74, 2, 'Coyame Del Sotol'),
    @(80, 2, 'Hidalgo Del Parral'),
    @(97, 2, 'San Juan De Sabinas'),
    @(105, 1, 'Ciudad De México'),
    @(130, 2, 'Pánuco De Coronado'),
    @(135, 2, 'San Juan De Guadalupe'),
    @(141, 1, 'Estado De México'),
    @(141, 2, 'Acambay De Ruíz Castañeda'),
    @(144, 2, 'Atizapán De Zaragoza'),
    @(149, 2, 'Coacalco De Berriozábal'),
    @(152, 2, 'Ecatepec De Morelos'),
    @(159, 2, 'Naucalpan De Juárez'),
    @(162, 2, 'San Antonio La Isla'),
    @(163, 2, 'San Felipe Del Progreso'),
    @(164, 2, 'San Martín De Las Pirámides'),
    @(172, 2, 'Tenango Del Valle'),
    @(175, 2, 'Tlalnepantla De Baz'),
    @(183, 2, 'San Miguel De Allende'),
    @(184, 2, 'Apaseo El Alto'),
    @(185, 2, 'Apaseo El Grande'),
    @(193, 2, 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @(197, 2, 'Jaral Del Progreso'),
    @(207, 2, 'San Diego De La Unión'),
    @(209, 2, 'San Francisco Del Rincón'),
    @(211, 2, 'San Luis De La Paz'),
    @(212, 2, 'Santa Cruz De Juventino Rosas'),
    @(213, 2, 'Silao De La Victoria'),
    @(218, 2, 'Valle De Santiago'),
    @(224, 2, 'Acapulco De Juárez'),
    @(225, 2, 'Ajuchitlán Del Progreso'),
    @(228, 2, 'Atenango Del Río'),
    @(229, 2, 'Atoyac De Álvarez'),
    @(230, 2, 'Ayutla De Los Libres'),
    @(232, 2, 'Chilapa De Álvarez'),
    @(233, 2, 'Chilpancingo De Los Bravo'),
    @(237, 2, 'Coyuca De Benítez'),
    @(238, 2, 'Coyuca De Catalán'),
    @(240, 2, 'Cuetzala Del Progreso'),
    @(241, 2, 'Cutzamala De Pinzón'),
    @(245, 2, 'Huitzuco De Los Figueroa'),
    @(246, 2, 'Iguala De La Independencia'),
    @(247, 2, 'Zihuatanejo De Azueta'),
    @(258, 2, 'Taxco De Alarcón'),
    @(260, 2, 'Técpan De Galeana'),
    @(265, 2, 'Tlapa De Comonfort'),
    @(271, 2, 'Atotonilco El Grande'),
    @(279, 2, 'Huejutla De Reyes'),
    @(282, 2, 'Jacala De Ledezma'),
    @(288, 2, 'Mixquiahuala De Juárez'),
    @(290, 2, 'Pachuca De Soto'),
    @(292, 2, 'Progreso De Obregón'),
    @(295, 2, 'Santiago De Anaya'),
    @(298, 2, 'Tenango De Doria'),
    @(300, 2, 'Tepeji Del Río De Ocampo'),
    @(306, 2, 'Tulancingo De Bravo'),
    @(310, 2, 'Atotonilco El Alto'),
    @(311, 2, 'Autlán De Navarro'),
    @(315, 2, 'Encarnación De Díaz'),
    @(317, 2, 'Huejuquilla El Alto'),
    @(319, 2, 'Jilotlán De Los Dolores'),
    @(321, 2, 'Lagos De Moreno'),
    @(324, 2, 'Ojuelos De Jalisco'),
    @(327, 2, 'Santa María Del Oro'),
    @(329, 2, 'Talpa De Allende'),
    @(330, 2, 'Tamazula De Gordiano'),
    @(332, 2, 'Tepatitlán De Morelos'),
    @(334, 2, 'Tizapán El Alto'),
    @(340, 2, 'Yahualica De González Gallo'),
    @(380, 2, 'Tiquicheo De Nicolás Romero'),
    @(405, 2, 'Santa María Del Oro'),
    @(418, 2, 'Mier Y Noriega'),
    @(422, 2, 'San Nicolás De Los Garza'),
    @(426, 2, 'Acatlán De Pérez Figueroa'),
    @(428, 2, 'Coicoyán De Las Flores'),
    @(430, 2, 'Guevea De Humboldt'),
    @(431, 2, 'Heroica Ciudad De Huajuapan De León'),
    @(432, 2, 'Heroica Ciudad De Tlaxiaco'),
    @(433, 2, 'Heroica Ciudad De Juchitán De Zaragoza'),
    @(436, 2, 'Miahuatlán De Porfirio Díaz'),
    @(437, 2, 'Oaxaca De Juárez'),
    @(471, 2, 'Santo Domingo De Morelos'),
    @(475, 2, 'Tamazulápam Del Espíritu Santo'),
    @(476, 2, 'Tataltepec De Valdés'),
    @(477, 2, 'Teotitlán De Flores Magón'),
    @(479, 2, 'Villa De Tututepec'),
    @(480, 2, 'Zapotitlán Del Río'),
    @(481, 2, 'Zimatlán De Álvarez'),
    @(487, 2, 'Ayotoxco De Guerrero'),
    @(500, 2, 'Izúcar De Matamoros'),
    @(506, 2, 'Los Reyes De Juárez'),
    @(509, 2, 'Palmar De Bravo'),
    @(514, 2, 'San Salvador El Seco'),
    @(516, 2, 'Tepexi De Rodríguez'),
    @(517, 2, 'Tepeyahualco De Cuauhtémoc'),
    @(525, 2, 'Tuzamapan De Galeana'),
    @(534, 2, 'Amealco De Bonfil'),
    @(538, 2, 'Jalpan De Serra'),
    @(539, 2, 'Landa De Matamoros'),
    @(541, 2, 'Pinal De Amoles'),
    @(543, 2, 'San Juan Del Río'),
    @(555, 2, 'Ciudad Del Maíz'),
    @(564, 2, 'Mexquitic De Carmona'),
    @(574, 2, 'Santa María Del Río'),
    @(582, 2, 'Tanquián De Escobedo'),
    @(586, 2, 'Villa De Arista'),
    @(587, 2, 'Villa De Guadalupe'),
    @(588, 2, 'Villa De Ramos'),
    @(589, 2, 'Villa De Reyes'),
    @(618, 2, 'Jalpa De Méndez'),
    @(654, 2, 'Nanacamilpa De Mariano Arista'),
    @(665, 2, 'Amatlán De Los Reyes'),
    @(670, 2, 'Camarón De Tejeda'),
    @(677, 2, 'Cosamaloapan De Carpio'),
    @(689, 2, 'Hueyapan De Ocampo'),
    @(695, 2, 'Juchique De Ferrer'),
    @(698, 2, 'Lerdo De Tejada'),
    @(700, 2, 'Martínez De La Torre'),
    @(711, 2, 'Poza Rica De Hidalgo'),
    @(715, 2, 'Sayula De Alemán'),
    @(717, 2, 'Soledad De Doblado'),
    @(735, 2, 'Vega De Alatorre'),
    @(739, 2, 'Zozocolco De Hidalgo'),
    @(760, 2, 'Nochistlán De Mejía'),
    @(767, 2, 'Teúl De González Ortega'),
    @(769, 2, 'Villa De Cos')
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $val = $u[2]
    $ws.Cells.Item($r, $c).Value = $val
}

# --- 3) Drop the footer / metadata rows (776-780) ---------------------
$ws.Range("A776:D780").ClearContents() | Out-Null
